$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4903499869063239
$ws.Range("C2").Value = 0.2576560982681002
$ws.Range("D2").Value = 0.07809297637275847
$ws.Range("E2").Value = 0.1454928687227106
$ws.Range("G2").Value = 0.8963065876550615
$ws.Range("H2").Value = 0.9556691858303168
$ws.Range("K2").Value = 0.2555851104106637
$ws.Range("L2").Value = 0.1911193299687284
$ws.Range("M2").Value = 0.1460821807584445
$ws.Range("O2").Value = 3.738499054475
$ws.Range("B3").Value = 0.4545619025943779
$ws.Range("C3").Value = 0.2580419132519474
$ws.Range("D3").Value = 0.0709065120186807
$ws.Range("E3").Value = 0.1461369765066323
$ws.Range("G3").Value = 0.8996451419604625
$ws.Range("H3").Value = 0.9612978588448655
$ws.Range("K3").Value = 0.2236868652359618
$ws.Range("L3").Value = 0.1885397940620379
$ws.Range("M3").Value = 0.1393518323464136
$ws.Range("O3").Value = 3.757073241339
$ws.Range("B4").Value = 0.4327238120229708
$ws.Range("C4").Value = 0.2583073982245736
$ws.Range("D4").Value = 0.06652725405521664
$ws.Range("E4").Value = 0.146588673436769
$ws.Range("G4").Value = 0.9021963064860472
$ws.Range("H4").Value = 0.9651256847939749
$ws.Range("K4").Value = 0.2040642024047941
$ws.Range("L4").Value = 0.1870465847903802
$ws.Range("M4").Value = 0.1352778560688925
$ws.Range("O4").Value = 3.770308741865875
$ws.Range("B5").Value = 0.4238593765453231
$ws.Range("C5").Value = 0.2584228052031676
$ws.Range("D5").Value = 0.06475106029182598
$ws.Range("E5").Value = 0.1467868990959857
$ws.Range("G5").Value = 0.9033619687822636
$ws.Range("H5").Value = 0.9667791329554589
$ws.Range("K5").Value = 0.1960589411263385
$ws.Range("L5").Value = 0.1864609441531329
$ws.Range("M5").Value = 0.1336324998310978
$ws.Range("O5").Value = 3.776162859814207
$ws.Range("B6").Value = 0.4223895605684049
$ws.Range("C6").Value = 0.2584424053743284
$ws.Range("D6").Value = 0.06445663213062858
$ws.Range("E6").Value = 0.1468206698749537
$ws.Range("G6").Value = 0.9035631385005019
$ws.Range("H6").Value = 0.9670593411887367
$ws.Range("K6").Value = 0.1947291510401925
$ws.Range("L6").Value = 0.1863650811694768
$ws.Range("M6").Value = 0.1333601883140219
$ws.Range("O6").Value = 3.777162751961669
$ws.Range("B7").Value = 0.4326041216173167
$ws.Range("C7").Value = 0.2583089253713347
$ws.Range("D7").Value = 0.06650326568227172
$ws.Range("E7").Value = 0.1465912894350048
$ws.Range("G7").Value = 0.9022115166885314
$ws.Range("H7").Value = 0.9651476047663721
$ws.Range("K7").Value = 0.2039562759686788
$ws.Range("L7").Value = 0.1870385940157249
$ws.Range("M7").Value = 0.1352556060510892
$ws.Range("O7").Value = 3.770385827620018
$ws.Range("B8").Value = 0.4779823779883827
$ws.Range("C8").Value = 0.2577832105118745
$ws.Range("D8").Value = 0.07560819367147076
$ws.Range("E8").Value = 0.1457033024433141
$ws.Range("G8").Value = 0.897353674209441
$ws.Range("H8").Value = 0.9575328370884506
$ws.Range("K8").Value = 0.2445945310606987
$ws.Range("L8").Value = 0.1902111310186498
$ws.Range("M8").Value = 0.1437494850730765
$ws.Range("O8").Value = 3.744523525642819
$ws.Range("B9").Value = 0.5680267870157536
$ws.Range("C9").Value = 0.2569778797254685
$ws.Range("D9").Value = 0.09372663687079807
$ws.Range("E9").Value = 0.1444071328390777
$ws.Range("G9").Value = 0.8918060547627391
$ws.Range("H9").Value = 0.9455467604486074
$ws.Range("K9").Value = 0.3239768735213886
$ws.Range("L9").Value = 0.1971493781396845
$ws.Range("M9").Value = 0.1608658837371451
$ws.Range("O9").Value = 3.70833008910617
$ws.Range("B10").Value = 0.6348062253255478
$ws.Range("C10").Value = 0.2565220548683911
$ws.Range("D10").Value = 0.1072001421842828
$ws.Range("E10").Value = 0.1437251620492646
$ws.Range("G10").Value = 0.8901584396550959
$ws.Range("H10").Value = 0.9385323085713679
$ws.Range("K10").Value = 0.3820954680558089
$ws.Range("L10").Value = 0.2026817121632973
$ws.Range("M10").Value = 0.1737175579576231
$ws.Range("O10").Value = 3.690588885755005
$ws.Range("B11").Value = 0.6653173498160641
$ws.Range("C11").Value = 0.2563438393434865
$ws.Range("D11").Value = 0.1133650525211038
$ws.Range("E11").Value = 0.1434733952395923
$ws.Range("G11").Value = 0.8899368298872474
$ws.Range("H11").Value = 0.9357293833357829
$ws.Range("K11").Value = 0.4084880565438596
$ws.Range("L11").Value = 0.2052924995102927
$ws.Range("M11").Value = 0.1796233162125134
$ws.Range("O11").Value = 3.684439246828646
$ws.Range("B12").Value = 0.6768897260279232
$ws.Range("C12").Value = 0.2562805157772843
$ws.Range("D12").Value = 0.1157046753496758
$ws.Range("E12").Value = 0.1433864458544178
$ws.Range("G12").Value = 0.8899288579367948
$ws.Range("H12").Value = 0.9347237036076308
$ws.Range("K12").Value = 0.4184752823021824
$ws.Range("L12").Value = 0.2062946157902843
$ws.Range("M12").Value = 0.1818681269954112
$ws.Range("O12").Value = 3.682386685620344
$ws.Range("B13").Value = 0.6743965959758498
$ws.Range("C13").Value = 0.2562939689376478
$ws.Range("D13").Value = 0.1152005692286906
$ws.Range("E13").Value = 0.1434047991326963
$ws.Range("G13").Value = 0.8899271965225921
$ws.Range("H13").Value = 0.9349378172696134
$ws.Range("K13").Value = 0.4163246768386557
$ws.Range("L13").Value = 0.2060781941535907
$ws.Range("M13").Value = 0.181384293934947
$ws.Range("O13").Value = 3.682816459361248
$ws.Range("B14").Value = 0.6662690494653987
$ws.Range("C14").Value = 0.2563385464049119
$ws.Range("D14").Value = 0.1135574325887774
$ws.Range("E14").Value = 0.1434660738313625
$ws.Range("G14").Value = 0.8899346517405746
$ws.Range("H14").Value = 0.9356455289706531
$ws.Range("K14").Value = 0.4093098561571082
$ws.Range("L14").Value = 0.2053746746711624
$ws.Range("M14").Value = 0.1798078299350863
$ws.Range("O14").Value = 3.684264846759959
$ws.Range("B15").Value = 0.6612930801062191
$ws.Range("C15").Value = 0.2563663926914188
$ws.Range("D15").Value = 0.1125516277302125
$ws.Range("E15").Value = 0.1435046983472255
$ws.Range("G15").Value = 0.8899491098443946
$ws.Range("H15").Value = 0.9360862782452273
$ws.Range("K15").Value = 0.4050121403449793
$ws.Range("L15").Value = 0.2049455006805516
$ws.Range("M15").Value = 0.1788432941807514
$ws.Range("O15").Value = 3.685187990226979
$ws.Range("B16").Value = 0.6328148709326626
$ws.Range("C16").Value = 0.2565342856954658
$ws.Range("D16").Value = 0.1067979671046686
$ws.Range("E16").Value = 0.143742790533782
$ws.Range("G16").Value = 0.8901835481134839
$ws.Range("H16").Value = 0.9387232880463188
$ws.Range("K16").Value = 0.3803696865006714
$ws.Range("L16").Value = 0.2025129791672811
$ws.Range("M16").Value = 0.173332788881126
$ws.Range("O16").Value = 3.691029428520039
$ws.Range("B17").Value = 0.6153779866173181
$ws.Range("C17").Value = 0.2566447266197045
$ws.Range("D17").Value = 0.1032774132095113
$ws.Range("E17").Value = 0.1439038139094198
$ws.Range("G17").Value = 0.8904626015488475
$ws.Range("H17").Value = 0.9404403342560954
$ws.Range("K17").Value = 0.3652402526200547
$ws.Range("L17").Value = 0.2010447621185421
$ws.Range("M17").Value = 0.1699674157826934
$ws.Range("O17").Value = 3.695104919053335
$ws.Range("B18").Value = 0.6053612911183563
$ws.Range("C18").Value = 0.2567109938901027
$ws.Range("D18").Value = 0.1012558503788483
$ws.Range("E18").Value = 0.1440019343848959
$ws.Range("G18").Value = 0.8906727919356143
$ws.Range("H18").Value = 0.941464456549383
$ws.Range("K18").Value = 0.3565339153127809
$ws.Range("L18").Value = 0.2002091435451661
$ws.Range("M18").Value = 0.1680373468673935
$ws.Range("O18").Value = 3.697629843041966
$ws.Range("B19").Value = 0.6019719827895074
$ws.Range("C19").Value = 0.2567339030816953
$ws.Range("D19").Value = 0.1005719632189965
$ws.Range("E19").Value = 0.1440361021721337
$ws.Range("G19").Value = 0.8907524913480813
$ws.Range("H19").Value = 0.9418174815443905
$ws.Range("K19").Value = 0.353585377471461
$ws.Range("L19").Value = 0.199927741288775
$ws.Range("M19").Value = 0.167384825054782
$ws.Range("O19").Value = 3.698515794747408
$ws.Range("B20").Value = 0.617232879489336
$ws.Range("C20").Value = 0.2566326861442718
$ws.Range("D20").Value = 0.1036518339966364
$ws.Range("E20").Value = 0.1438861031739975
$ws.Range("G20").Value = 0.8904277532336664
$ws.Range("H20").Value = 0.9402537722566677
$ws.Range("K20").Value = 0.3668512538743869
$ws.Range("L20").Value = 0.2012001397079786
$ws.Range("M20").Value = 0.1703250863224568
$ws.Range("O20").Value = 3.694652363385387
$ws.Range("B21").Value = 0.6686558104383948
$ws.Range("C21").Value = 0.256325340164409
$ws.Range("D21").Value = 0.1140399233834017
$ws.Range("E21").Value = 0.1434478484286714
$ws.Range("G21").Value = 0.8899304005060458
$ws.Range("H21").Value = 0.9354361449796897
$ws.Range("K21").Value = 0.4113704748671694
$ws.Range("L21").Value = 0.2055809503941646
$ws.Range("M21").Value = 0.1802706478639493
$ws.Range("O21").Value = 3.683831925426688
$ws.Range("B22").Value = 0.7023710692662064
$ws.Range("C22").Value = 0.2561487230576986
$ws.Range("D22").Value = 0.120858879778595
$ws.Range("E22").Value = 0.143210314201955
$ws.Range("G22").Value = 0.890048044134474
$ws.Range("H22").Value = 0.9326123383053471
$ws.Range("K22").Value = 0.4404249001865139
$ws.Range("L22").Value = 0.2085225342316193
$ws.Range("M22").Value = 0.1868197168587074
$ws.Range("O22").Value = 3.678369869371949
$ws.Range("B23").Value = 0.6843669798647625
$ws.Range("C23").Value = 0.2562407770172115
$ws.Range("D23").Value = 0.117216765046507
$ws.Range("E23").Value = 0.1433326231338619
$ws.Range("G23").Value = 0.8899447375599152
$ws.Range("H23").Value = 0.9340897598573719
$ws.Range("K23").Value = 0.4249219587647417
$ws.Range("L23").Value = 0.2069453966024639
$ws.Range("M23").Value = 0.1833199062262665
$ws.Range("O23").Value = 3.681137799257897
$ws.Range("B24").Value = 0.6163942574521002
$ws.Range("C24").Value = 0.2566381210031068
$ws.Range("D24").Value = 0.1034825507243085
$ws.Range("E24").Value = 0.1438940929211885
$ws.Range("G24").Value = 0.8904433531636755
$ws.Range("H24").Value = 0.9403380017725169
$ws.Range("K24").Value = 0.366122945789698
$ws.Range("L24").Value = 0.2011298670848873
$ws.Range("M24").Value = 0.1701633687325597
$ws.Range("O24").Value = 3.694856397280518
$ws.Range("B25").Value = 0.5435563263725953
$ws.Range("C25").Value = 0.2571717661713961
$ws.Range("D25").Value = 0.08879676740340869
$ws.Range("E25").Value = 0.1447102324416036
$ws.Range("G25").Value = 0.8928806088690777
$ws.Range("H25").Value = 0.9484743418894794
$ws.Range("K25").Value = 0.3025365151589767
$ws.Range("L25").Value = 0.1951958799654108
$ws.Range("M25").Value = 0.1561866333289466
$ws.Range("O25").Value = 3.716566968525228
